# Updates scraped market-board price/profit figures across the "Ultima
# Profits" leve-crafting workbook. Each worksheet corresponds to a crafting
# job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR); columns H:N hold the
# currentAveragePrice* / LevePrice* / LeveProfit* figures that the scheduled
# runner refreshes from the market API.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H40").Value = 1582.174
$ws.Range("I40").Value = 1720
$ws.Range("J40").Value = 1476.1538
$ws.Range("K40").Value = 1720
$ws.Range("L40").Value = 1476.1538
$ws.Range("M40").Value = -1545
$ws.Range("N40").Value = -1826.1538

$ws.Range("H62").Value = 1196.5
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 1262
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 1262
$ws.Range("M62").Value = -376
$ws.Range("N62").Value = -2510

$ws.Range("H65").Value = 1196.5
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 1262
$ws.Range("K65").Value = 1000
$ws.Range("L65").Value = 6310
$ws.Range("M65").Value = -1880
$ws.Range("N65").Value = -12550

$ws.Range("H100").Value = 2006
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H106").Value = 457281.38
$ws.Range("I106").Value = 502559.5
$ws.Range("K106").Value = 502559.5
$ws.Range("M106").Value = -501928.5

$ws.Range("H116").Value = 2998.4666
$ws.Range("J116").Value = 3215.3
$ws.Range("L116").Value = 3215.3
$ws.Range("N116").Value = -10099.3

$ws.Range("H137").Value = 10000771
$ws.Range("I137").Value = 695.2353000000001
$ws.Range("J137").Value = 66667868
$ws.Range("K137").Value = 2085.7059
$ws.Range("L137").Value = 200003604
$ws.Range("M137").Value = 464.2941000000001
$ws.Range("N137").Value = -200008704

$ws.Range("H138").Value = 6174401.5
$ws.Range("I138").Value = 8773164
$ws.Range("J138").Value = 2339.625
$ws.Range("K138").Value = 26319492
$ws.Range("L138").Value = 7018.875
$ws.Range("M138").Value = -26314352
$ws.Range("N138").Value = -17298.875

# ---------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 11284.333
$ws.Range("I32").Value = 11094.4
$ws.Range("J32").Value = 12234
$ws.Range("K32").Value = 11094.4
$ws.Range("L32").Value = 12234
$ws.Range("M32").Value = -10807.4
$ws.Range("N32").Value = -12808

$ws.Range("H61").Value = 12501508
$ws.Range("I61").Value = 14707392
$ws.Range("J61").Value = 1498
$ws.Range("K61").Value = 14707392
$ws.Range("L61").Value = 1498
$ws.Range("M61").Value = -14707180
$ws.Range("N61").Value = -1922

$ws.Range("H136").Value = 12501508
$ws.Range("I136").Value = 14707392
$ws.Range("J136").Value = 1498
$ws.Range("K136").Value = 44122176
$ws.Range("L136").Value = 4494
$ws.Range("M136").Value = -44119626
$ws.Range("N136").Value = -9594

# ---------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H10").Value = 5006
$ws.Range("J10").Value = 5006
$ws.Range("L10").Value = 5006
$ws.Range("N10").Value = -5286

$ws.Range("H94").Value = 801.7059
$ws.Range("I94").Value = 680.62067
$ws.Range("J94").Value = 1504
$ws.Range("K94").Value = 680.62067
$ws.Range("L94").Value = 1504
$ws.Range("M94").Value = -229.62067
$ws.Range("N94").Value = -2406

# ---------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H58").Value = 1306.8286
$ws.Range("I58").Value = 605.7917
$ws.Range("J58").Value = 2836.3635
$ws.Range("K58").Value = 605.7917
$ws.Range("L58").Value = 2836.3635
$ws.Range("M58").Value = -402.7917
$ws.Range("N58").Value = -3242.3635

$ws.Range("H136").Value = 1306.8286
$ws.Range("I136").Value = 605.7917
$ws.Range("J136").Value = 2836.3635
$ws.Range("K136").Value = 1817.3751
$ws.Range("L136").Value = 8509.0905
$ws.Range("M136").Value = 732.6249
$ws.Range("N136").Value = -13609.0905

# ---------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H132").Value = 3906.6667
$ws.Range("I132").Value = 860
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 7740
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -5210
$ws.Range("N132").Value = -95060

# ---------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H126").Value = 4148.5835
$ws.Range("I126").Value = 2495.75
$ws.Range("J126").Value = 4975
$ws.Range("K126").Value = 7487.25
$ws.Range("L126").Value = 14925
$ws.Range("M126").Value = -5017.25
$ws.Range("N126").Value = -19865

$ws.Range("H132").Value = 4023.15
$ws.Range("I132").Value = 2971.5715
$ws.Range("J132").Value = 6476.8335
$ws.Range("K132").Value = 8914.7145
$ws.Range("L132").Value = 19430.5005
$ws.Range("M132").Value = -6384.7145
$ws.Range("N132").Value = -24490.5005

# ---------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H16").Value = 1602.5
$ws.Range("I16").Value = 1930.5
$ws.Range("J16").Value = 946.5
$ws.Range("K16").Value = 1930.5
$ws.Range("L16").Value = 946.5
$ws.Range("M16").Value = -1760.5
$ws.Range("N16").Value = -1286.5

$ws.Range("H40").Value = 18000
$ws.Range("I40").Value = 26500
$ws.Range("K40").Value = 26500
$ws.Range("M40").Value = -26364

$ws.Range("H88").Value = 22057
$ws.Range("I88").Value = 10171
$ws.Range("J88").Value = 28000
$ws.Range("K88").Value = 10171
$ws.Range("L88").Value = 28000
$ws.Range("M88").Value = -9743
$ws.Range("N88").Value = -28856

$ws.Range("H91").Value = 22057
$ws.Range("I91").Value = 10171
$ws.Range("J91").Value = 28000
$ws.Range("K91").Value = 10171
$ws.Range("L91").Value = 28000
$ws.Range("M91").Value = -8689
$ws.Range("N91").Value = -30964

$ws.Range("H100").Value = 2910
$ws.Range("I100").Value = 3075
$ws.Range("J100").Value = 2745
$ws.Range("K100").Value = 3075
$ws.Range("L100").Value = 2745
$ws.Range("M100").Value = -2534
$ws.Range("N100").Value = -3827

# ---------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H5").Value = 2858671.5
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 3335066.8
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 3335066.8
$ws.Range("M5").Value = -188
$ws.Range("N5").Value = -3335290.8

$ws.Range("H62").Value = 11628.714
$ws.Range("I62").Value = 6781.125
$ws.Range("J62").Value = 18092.166
$ws.Range("K62").Value = 6781.125
$ws.Range("L62").Value = 18092.166
$ws.Range("M62").Value = -6157.125
$ws.Range("N62").Value = -19340.166

$ws.Range("H65").Value = 11628.714
$ws.Range("I65").Value = 6781.125
$ws.Range("J65").Value = 18092.166
$ws.Range("K65").Value = 33905.625
$ws.Range("L65").Value = 90460.83
$ws.Range("M65").Value = -30785.625
$ws.Range("N65").Value = -96700.83
